$wb = $excel.ActiveWorkbook

# --- 1) Duplicate the "Add Devices" sheet, place the copy right after it,
#        and rename the copy to "Sheet1" (new tab added by this commit). ---
$src = $wb.Worksheets.Item("Add Devices")
$srcIndex = $src.Index
$src.Copy($null, $src)
$newSheet = $wb.Worksheets.Item($srcIndex + 1)
$newSheet.Name = "Sheet1"
$newSheet.Range("C7").Select()

# --- 2) Re-activate the original sheet and update its row 10 data/values. ---
$ws = $wb.Worksheets.Item("Add Devices")
$ws.Activate()

$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("H10").Value = 1
$ws.Range("J10").Value = 1
$ws.Range("K10").Value = "Yes"
$ws.Range("L10").Value = $false
$ws.Range("N10").Value = "Other Slot Cards  (6"
$ws.Range("O10").Value = "NA"

# Row 10 no longer needs the taller, custom row height.
$ws.Rows.Item(10).AutoFit()

# Update the active selection shown on the "Add Devices" tab.
$ws.Range("N6").Select()
